$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 4387.2,              4741, 3890, 0.05135540167490641),
    @(1, 4189.9,              4512, 3577, 0.05696694056193034),
    @(2, 4850.6,              5145, 4607, 0.05506304105122884),
    @(3, 5112.033333333334,   5522, 4570, 0.05476688543955485),
    @(4, 4328.2,              4562, 3976, 0.0528297742207845),
    @(5, 4986,                5314, 4421, 0.05610188643137614),
    @(6, 4425.866666666667,   4815, 3915, 0.05510924657185872),
    @(7, 4518.7,              4878, 4143, 0.05585019588470459),
    @(8, 4680.5,              5155, 4096, 0.05555754502614339),
    @(9, 4780.966666666666,   5243, 4389, 0.05167122681935628)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
